$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3151.3333
$ws.Range("I76").Value = 3199
$ws.Range("J76").Value = 3127.5
$ws.Range("K76").Value = 3199
$ws.Range("L76").Value = 3127.5
$ws.Range("M76").Value = -2884
$ws.Range("N76").Value = -3757.5

$ws.Range("H79").Value = 3151.3333
$ws.Range("I79").Value = 3199
$ws.Range("J79").Value = 3127.5
$ws.Range("K79").Value = 3199
$ws.Range("L79").Value = 3127.5
$ws.Range("M79").Value = -2107
$ws.Range("N79").Value = -5311.5

$ws.Range("H86").Value = 3425
$ws.Range("I86").Value = 3425
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 3425
$ws.Range("L86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -2302

$ws.Range("H89").Value = 3425
$ws.Range("I89").Value = 3425
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 17125
$ws.Range("L89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -11509

$ws.Range("H112").Value = 1836.8918
$ws.Range("I112").Value = 500
$ws.Range("J112").Value = 1874.0278
$ws.Range("K112").Value = 1500
$ws.Range("L112").Value = 5622.0834
$ws.Range("M112").Value = -392
$ws.Range("N112").Value = -7838.0834

$ws.Range("H113").Value = 74328.63
$ws.Range("I113").Value = 115258.43
$ws.Range("J113").Value = 2701.5
$ws.Range("K113").Value = 115258.43
$ws.Range("L113").Value = 2701.5
$ws.Range("M113").Value = -112004.43
$ws.Range("N113").Value = -9209.5

$ws.Range("H138").Value = 3806.6296
$ws.Range("I138").Value = 11269.75
$ws.Range("J138").Value = 2508.6956
$ws.Range("K138").Value = 33809.25
$ws.Range("L138").Value = 7526.0868
$ws.Range("M138").Value = -28669.25
$ws.Range("N138").Value = -17806.0868

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1848.6428
$ws.Range("I61").Value = 1698.1111
$ws.Range("J61").Value = 2119.6
$ws.Range("K61").Value = 1698.1111
$ws.Range("L61").Value = 2119.6
$ws.Range("M61").Value = -1486.1111
$ws.Range("N61").Value = -2543.6

$ws.Range("H74").Value = 1024.2195
$ws.Range("I74").Value = 535.3823
$ws.Range("J74").Value = 3398.5715
$ws.Range("K74").Value = 535.3823
$ws.Range("L74").Value = 3398.5715
$ws.Range("M74").Value = 338.6177
$ws.Range("N74").Value = -5146.5715

$ws.Range("H77").Value = 1024.2195
$ws.Range("I77").Value = 535.3823
$ws.Range("J77").Value = 3398.5715
$ws.Range("K77").Value = 2676.9115
$ws.Range("L77").Value = 16992.8575
$ws.Range("M77").Value = 1691.0885
$ws.Range("N77").Value = -25728.8575

$ws.Range("H122").Value = 2158.6924
$ws.Range("I122").Value = 2159
$ws.Range("J122").Value = 2157
$ws.Range("K122").Value = 6477
$ws.Range("L122").Value = 6471
$ws.Range("M122").Value = -4027
$ws.Range("N122").Value = -11371

$ws.Range("H132").Value = 1407.3684
$ws.Range("I132").Value = 1373.7222
$ws.Range("J132").Value = 2013
$ws.Range("K132").Value = 4121.1666
$ws.Range("L132").Value = 6039
$ws.Range("M132").Value = -1591.1666
$ws.Range("N132").Value = -11099

$ws.Range("H136").Value = 1848.6428
$ws.Range("I136").Value = 1698.1111
$ws.Range("J136").Value = 2119.6
$ws.Range("K136").Value = 5094.3333
$ws.Range("L136").Value = 6358.799999999999
$ws.Range("M136").Value = -2544.3333
$ws.Range("N136").Value = -11458.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1003277.56
$ws.Range("I86").Value = 1115181.4
$ws.Range("J86").Value = 667566.3
$ws.Range("K86").Value = 1115181.4
$ws.Range("L86").Value = 667566.3
$ws.Range("M86").Value = -1114058.4
$ws.Range("N86").Value = -669812.3

$ws.Range("H89").Value = 1003277.56
$ws.Range("I89").Value = 1115181.4
$ws.Range("J89").Value = 667566.3
$ws.Range("K89").Value = 5575907
$ws.Range("L89").Value = 3337831.5
$ws.Range("M89").Value = -5570291
$ws.Range("N89").Value = -3349063.5

$ws.Range("H134").Value = 6050.5186
$ws.Range("I134").Value = 6938.864
$ws.Range("J134").Value = 2141.8
$ws.Range("K134").Value = 20816.592
$ws.Range("L134").Value = 6425.400000000001
$ws.Range("M134").Value = -18281.592
$ws.Range("N134").Value = -11495.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2740.6667
$ws.Range("I31").Value = 2583.9
$ws.Range("J31").Value = 2936.625
$ws.Range("K31").Value = 2583.9
$ws.Range("L31").Value = 2936.625
$ws.Range("M31").Value = -2288.9
$ws.Range("N31").Value = -3526.625

$ws.Range("H34").Value = 2740.6667
$ws.Range("I34").Value = 2583.9
$ws.Range("J34").Value = 2936.625
$ws.Range("K34").Value = 2583.9
$ws.Range("L34").Value = 2936.625
$ws.Range("M34").Value = -2381.9
$ws.Range("N34").Value = -3340.625

$ws.Range("H99").Value = 2958.111
$ws.Range("I99").Value = 2090.5715
$ws.Range("J99").Value = 5994.5
$ws.Range("K99").Value = 2090.5715
$ws.Range("L99").Value = 5994.5
$ws.Range("M99").Value = -592.5715
$ws.Range("N99").Value = -8990.5

$ws.Range("H126").Value = 2958.111
$ws.Range("I126").Value = 2090.5715
$ws.Range("J126").Value = 5994.5
$ws.Range("K126").Value = 6271.7145
$ws.Range("L126").Value = 17983.5
$ws.Range("M126").Value = -3801.7145
$ws.Range("N126").Value = -22923.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 6655.4116
$ws.Range("I113").Value = 34167
$ws.Range("J113").Value = 760.0714
$ws.Range("K113").Value = 102501
$ws.Range("L113").Value = 2280.2142
$ws.Range("M113").Value = -100331
$ws.Range("N113").Value = -6620.2142

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("M80").ClearContents()

$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("M83").ClearContents()

$ws.Range("H106").Value = 25000
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 25000
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 25000
$ws.Range("N106").Value = -27524

$ws.Range("H122").Value = 1732.25
$ws.Range("I122").Value = 1618.5714
$ws.Range("J122").Value = 1997.5
$ws.Range("K122").Value = 4855.7142
$ws.Range("L122").Value = 5992.5
$ws.Range("M122").Value = -2405.7142
$ws.Range("N122").Value = -10892.5

$ws.Range("H132").Value = 1482362.2
$ws.Range("I132").Value = 2139273
$ws.Range("J132").Value = 4313.125
$ws.Range("K132").Value = 6417819
$ws.Range("L132").Value = 12939.375
$ws.Range("M132").Value = -6415289
$ws.Range("N132").Value = -17999.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2837.25
$ws.Range("I7").Value = 2618.1428
$ws.Range("J7").Value = 3144
$ws.Range("K7").Value = 2618.1428
$ws.Range("L7").Value = 3144
$ws.Range("M7").Value = -2506.1428
$ws.Range("N7").Value = -3368

$ws.Range("H22").Value = 2069.3333
$ws.Range("I22").Value = 2320.2856
$ws.Range("J22").Value = 1718
$ws.Range("K22").Value = 2320.2856
$ws.Range("L22").Value = 1718
$ws.Range("M22").Value = -2025.2856
$ws.Range("N22").Value = -2308

$ws.Range("H27").Value = 2069.3333
$ws.Range("I27").Value = 2320.2856
$ws.Range("J27").Value = 1718
$ws.Range("K27").Value = 2320.2856
$ws.Range("L27").Value = 1718
$ws.Range("M27").Value = -2213.2856
$ws.Range("N27").Value = -1932

$ws.Range("H55").Value = 504.09525
$ws.Range("I55").Value = 484.9091
$ws.Range("J55").Value = 525.2
$ws.Range("K55").Value = 484.9091
$ws.Range("L55").Value = 525.2
$ws.Range("M55").Value = -311.9091
$ws.Range("N55").Value = -871.2

$ws.Range("H82").Value = 1960.6666
$ws.Range("I82").Value = 1382.2858
$ws.Range("J82").Value = 3985
$ws.Range("K82").Value = 1382.2858
$ws.Range("L82").Value = 3985
$ws.Range("M82").Value = -1021.2858
$ws.Range("N82").Value = -4707

$ws.Range("H85").Value = 1960.6666
$ws.Range("I85").Value = 1382.2858
$ws.Range("J85").Value = 3985
$ws.Range("K85").Value = 1382.2858
$ws.Range("L85").Value = 3985
$ws.Range("M85").Value = -134.2858000000001
$ws.Range("N85").Value = -6481

$ws.Range("H93").Value = 15152329
$ws.Range("I93").Value = 827.8461
$ws.Range("J93").Value = 37037830
$ws.Range("K93").Value = 827.8461
$ws.Range("L93").Value = 37037830
$ws.Range("M93").Value = 420.1539
$ws.Range("N93").Value = -37040326

$ws.Range("H126").Value = 2837.25
$ws.Range("I126").Value = 2618.1428
$ws.Range("J126").Value = 3144
$ws.Range("K126").Value = 7854.428400000001
$ws.Range("L126").Value = 9432
$ws.Range("M126").Value = -5384.428400000001
$ws.Range("N126").Value = -14372

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 135323.27
$ws.Range("I122").Value = 97969.5
$ws.Range("J122").Value = 234933.33
$ws.Range("K122").Value = 293908.5
$ws.Range("L122").Value = 704799.99
$ws.Range("M122").Value = -291458.5
$ws.Range("N122").Value = -709699.99

$ws.Range("H136").Value = 15433825
$ws.Range("I136").Value = 25254368
$ws.Range("J136").Value = 1542.0714
$ws.Range("K136").Value = 75763104
$ws.Range("L136").Value = 4626.2142
$ws.Range("M136").Value = -75760554
$ws.Range("N136").Value = -9726.2142
